# Updates the crypto price/volume table (columns D = Price, E = Volume(1h))
# to the latest scraped snapshot. Values that could be mis-parsed by Excel's
# automatic type detection as plain numbers (losing significant trailing
# zeros, e.g. "1.000" -> 1, or flipping to scientific notation, e.g.
# "0.000007964" -> 7.964E-06) are temporarily forced to Text format, written,
# and then restored to their original cell style so no formatting changes
# are left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.632.28'
$ws.Range("E2").Value = '  -5.90%  '
$ws.Range("D3").Value = '1.805.15'
$ws.Range("E3").Value = '  -5.20%  '
$ws.Range("E4").Value = '  +0.13%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '275.51'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -10.01%  '
$ws.Range("E6").Value = '  +0.14%  '
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5058'
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = '  -6.52%  '
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3501'
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  -8.13%  '
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.75'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  -4.86%  '
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06617'
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  -9.19%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.97'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -9.89%  '
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8342'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  -7.65%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07758'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  -5.19%  '
$ws.Range("D14").Value = '1.801.99'
$ws.Range("E14").Value = '  +54.02%  '
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.074'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  -5.20%  '
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.51'
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  -8.68%  '
$ws.Range("E17").Value = '  +0.13%  '
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.90'
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = '  -6.69%  '
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = $style
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007964'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -7.90%  '
$ws.Range("D21").Value = '25.691.10'
$ws.Range("E21").Value = '  -5.78%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.721'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  -6.52%  '
$ws.Range("D23").Value = '2.035.07'
$ws.Range("E23").Value = '  +55.40%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.03'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -7.27%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.053'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -7.15%  '
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.49'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  -4.06%  '
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.107'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  -8.80%  '
$ws.Range("E28").Value = '  -5.71%  '
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.93'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -7.88%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '108.24'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -7.46%  '
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.317'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  -11.20%  '
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.197'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  -10.23%  '
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08789'
$ws.Range("D33").Style = $style
$ws.Range("E34").Value = '  -5.44%  '
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7224'
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  -13.31%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.126'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  -8.00%  '
$ws.Range("E37").Value = '  -5.10%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9996'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  +0.14%  '
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.024'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  -8.98%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01861'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -7.18%  '
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5174'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -13.39%  '
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.275'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  -15.63%  '
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9546'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  -12.00%  '
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '114.60'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  -1.51%  '
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.171'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  -7.31%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.996'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -13.82%  '
$ws.Range("E47").Value = '  +0.10%  '
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4565'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -10.97%  '
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1376'
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  -10.05%  '
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.297'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  -8.89%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.87'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -6.24%  '
